# Applies the "small edits to slides and paper" commit:
# text tweaks scattered across slides 2-3 (The Question / The Data),
# 4-5 (The [Modeling] Approach[, continued]), 6 (Challenges),
# 7 (Results) and 8 (Next Steps).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# Slide 2 - "The Question"
# ---------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$s2Body = $s2.Shapes.Item(2).TextFrame.TextRange

$s2Body.Paragraphs(1,1).Text = "Can visual features be used for predictive modeling of security prices? "
$s2Body.Paragraphs(2,1).Text = "“Technical analysis” attempts to do so, with humans as the “algorithm”… but is basically bunk"

# ---------------------------------------------------------------
# Slide 3 - "The Data"
# ---------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$s3Body = $s3.Shapes.Item(2).TextFrame.TextRange

$s3Body.Paragraphs(5,1).Text = "Time indexing is important, to keep our models honest by barring them from peering into the future"

# ---------------------------------------------------------------
# Slide 4 - "The Approach" -> "The Modeling Approach"
# ---------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "The Modeling Approach"

$s4Body = $s4.Shapes.Item(2).TextFrame.TextRange
$s4Body.Paragraphs(3,1).Text = "Window length is somewhat arbitrary, but impacts predictions, so can be tuned depending on specific use case (macro vs. HFT, e.g.)"
$s4Body.Paragraphs(5,1).Text = "Use computer vision (CV) to extract features:"
$s4Body.Paragraphs(6,1).Text = "Directly from images using linear feature extraction"
$s4Body.Paragraphs(7,1).Text = "Indirectly, by converting the graph to pixel intensity data, then unrolling the pixel data matrix into a vector 1000s of new features for each observation, “generated” from the underlying price data"

# ---------------------------------------------------------------
# Slide 5 - "Approach, continued" -> "Modeling Approach, continued"
# ---------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "Modeling Approach, continued"

# ---------------------------------------------------------------
# Slide 6 - "Challenges"
# ---------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$s6Body = $s6.Shapes.Item(2).TextFrame.TextRange

$s6Body.Paragraphs(2,1).Text = "Having patience while trying to train models on massive datasets (~ 6,000 x 30,000 matrix of raw pixel data)"
$s6Body.Paragraphs(3,1).Text = "Not being able to use regular TTS / Cross-Validation due to time series data"

# ---------------------------------------------------------------
# Slide 7 - "Results"
# ---------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$s7Body = $s7.Shapes.Item(2).TextFrame.TextRange

$s7Body.Paragraphs(2,1).Text = "Models varied in performance over these data; high end performance was impressive"

# ---------------------------------------------------------------
# Slide 8 - "Next Steps"
# ---------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$s8Body = $s8.Shapes.Item(2).TextFrame.TextRange

$s8Body.Paragraphs(5,1).Text = "Additional time-series features, e.g. weighted trailing feature stats baked into current observation"
$s8Body.Paragraphs(7,1).Text = "Further “productionizing” of analytic processes as Python is very nice for this sort of thing:"
